# Add three new groups (Gruppe 19-21) with their album cover paths and
# initial score of 0, continuing the existing table pattern on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: group names (fills the two previously-blank rows 20-21 and adds row 22)
$ws.Range("A20").Value = "Gruppe 19"
$ws.Range("A21").Value = "Gruppe 20"
$ws.Range("A22").Value = "Gruppe 21"

# Column B: matching album cover image paths
$ws.Range("B20").Value = "/album_covers/artist_19.png"
$ws.Range("B21").Value = "/album_covers/artist_20.png"
$ws.Range("B22").Value = "/album_covers/artist_21.png"

# Column C: starting score
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("C22").Value = 0

# Match the formatting used by the rest of the data rows (copy row 19's style)
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Mirror the cursor/view position left behind by the edit
$ws.Range("E22").Select() | Out-Null
